$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The first 16 data rows (old rows 2-17, corresponding to dates 30864..32234)
# are obsolete/typo'd and must be removed. Deleting these rows shifts all
# subsequent rows up by 16, so what was row 18 (date 32325) becomes row 2, etc.
$ws.Range("A2:B17").EntireRow.Delete() | Out-Null
